# "contingencies with rene fine"
# Adds a small 2x2 block to Sheet1:
#   B1 = 0                      (bold, thin box border, centered/top aligned)
#   A2 = 0                      (same style as B1)
#   B2 = "disconnected_elements" (plain shared-string text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- formatting for B1 -------------------------------------------------
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment = -4160     # xlTop
$r1.Borders.LineStyle = 1         # xlContinuous
$r1.Borders.Weight = 2            # xlThin

# --- propagate the identical format to A2 ------------------------------
# (copy/paste-format keeps both cells pointed at the same cellXfs entry,
# instead of each Range(...) call minting its own style record)
$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial(-4122)           # xlPasteFormats
